$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Row 2: ID column was stored as text "1" - change to a real number 1
$ws.Range("A2").Value = 1

# numeroCuenta (column N) value correction on rows 2 and 3
$ws.Range("N2").Value = "406-101480-06"
$ws.Range("N3").Value = "406-101480-06"

# Updated expected result codes
$ws.Range("P2").Value = 798514
$ws.Range("Q3").Value = 798515
$ws.Range("R3").Value = 798516
